$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: fill in G19/H19 values, and add new I19/J19 cells with values,
# copying number formatting (fill/border) from the analogous row 16 cells
# so that the new cells pick up the same styles used elsewhere in the sheet.
$ws.Range("G19").Value = 5
$ws.Range("H19").Value = 5

$ws.Range("I16").Copy()
$ws.Range("I19").PasteSpecial(-4122)
$ws.Range("I19").Value = 5

$ws.Range("J16").Copy()
$ws.Range("J19").PasteSpecial(-4122)
$ws.Range("J19").Value = 5

# Row 28: add new J28 cell (default/no special style) with value 5.
$ws.Range("J28").Value = 5

# Update the active selection in the bottom-right frozen pane to J19.
$ws.Range("J19").Select()
